$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for trial 1 (row 2) and trial 2 (row 3) had their
# condition values (columns B:J) swapped between the two trials, while
# the trial index in column A stays the same.

$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 11
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = "instr_dim1_1"

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 16
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = "instr_dim1_1"
